# Auto-generated COM-interop edit script.
# Updates cached market-price / profit figures (columns H-N) on several
# Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1286.762
$ws.Range("I40").Value = 1388.1904
$ws.Range("J40").Value = 1185.3334
$ws.Range("K40").Value = 1388.1904
$ws.Range("L40").Value = 1185.3334
$ws.Range("M40").Value = -1213.1904
$ws.Range("N40").Value = -1535.3334
$ws.Range("H43").Value = 1090.9375
$ws.Range("I43").Value = 829.1667
$ws.Range("J43").Value = 1248
$ws.Range("K43").Value = 829.1667
$ws.Range("L43").Value = 1248
$ws.Range("M43").Value = -760.1667
$ws.Range("N43").Value = -1386
$ws.Range("H53").Value = 231.38461
$ws.Range("I53").Value = 129.8
$ws.Range("J53").Value = 294.875
$ws.Range("K53").Value = 129.8
$ws.Range("L53").Value = 294.875
$ws.Range("M53").Value = 507.2
$ws.Range("N53").Value = -1568.875
$ws.Range("H98").Value = 28328.379
$ws.Range("I98").Value = 44390.957
$ws.Range("J98").Value = 1939.8572
$ws.Range("K98").Value = 44390.957
$ws.Range("L98").Value = 1939.8572
$ws.Range("M98").Value = -42892.957
$ws.Range("N98").Value = -4935.8572
$ws.Range("H122").Value = 28328.379
$ws.Range("I122").Value = 44390.957
$ws.Range("J122").Value = 1939.8572
$ws.Range("K122").Value = 133172.871
$ws.Range("L122").Value = 5819.571599999999
$ws.Range("M122").Value = -130722.871
$ws.Range("N122").Value = -10719.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10193.167
$ws.Range("I32").Value = 5482.1665
$ws.Range("K32").Value = 5482.1665
$ws.Range("M32").Value = -5195.1665
$ws.Range("H61").Value = 1628.7222
$ws.Range("I61").Value = 1408.5834
$ws.Range("J61").Value = 2069
$ws.Range("K61").Value = 1408.5834
$ws.Range("L61").Value = 2069
$ws.Range("M61").Value = -1196.5834
$ws.Range("N61").Value = -2493
$ws.Range("H74").Value = 45456720
$ws.Range("I74").Value = 71430380
$ws.Range("J74").Value = 2800
$ws.Range("K74").Value = 71430380
$ws.Range("L74").Value = 2800
$ws.Range("M74").Value = -71429506
$ws.Range("N74").Value = -4548
$ws.Range("H77").Value = 45456720
$ws.Range("I77").Value = 71430380
$ws.Range("J77").Value = 2800
$ws.Range("K77").Value = 357151900
$ws.Range("L77").Value = 14000
$ws.Range("M77").Value = -357147532
$ws.Range("N77").Value = -22736
$ws.Range("H132").Value = 2063.3914
$ws.Range("I132").Value = 1463.3226
$ws.Range("J132").Value = 3303.5334
$ws.Range("K132").Value = 4389.9678
$ws.Range("L132").Value = 9910.600199999999
$ws.Range("M132").Value = -1859.9678
$ws.Range("N132").Value = -14970.6002
$ws.Range("H136").Value = 1628.7222
$ws.Range("I136").Value = 1408.5834
$ws.Range("J136").Value = 2069
$ws.Range("K136").Value = 4225.7502
$ws.Range("L136").Value = 6207
$ws.Range("M136").Value = -1675.7502
$ws.Range("N136").Value = -11307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5842.857
$ws.Range("I7").Value = 700
$ws.Range("J7").Value = 7900
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 7900
$ws.Range("M7").Value = -587
$ws.Range("N7").Value = -8126
$ws.Range("H22").Value = 399.7857
$ws.Range("I22").Value = 358.9091
$ws.Range("J22").Value = 549.6667
$ws.Range("K22").Value = 358.9091
$ws.Range("L22").Value = 549.6667
$ws.Range("M22").Value = -185.9091
$ws.Range("N22").Value = -895.6667
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
$ws.Range("H134").Value = 11906368
$ws.Range("I134").Value = 14707230
$ws.Range("J134").Value = 2703.5
$ws.Range("K134").Value = 44121690
$ws.Range("L134").Value = 8110.5
$ws.Range("M134").Value = -44119155
$ws.Range("N134").Value = -13180.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 36381.21
$ws.Range("I68").Value = 20089.334
$ws.Range("J68").Value = 39435.938
$ws.Range("K68").Value = 20089.334
$ws.Range("L68").Value = 39435.938
$ws.Range("M68").Value = -19340.334
$ws.Range("N68").Value = -40933.938
$ws.Range("H71").Value = 36381.21
$ws.Range("I71").Value = 20089.334
$ws.Range("J71").Value = 39435.938
$ws.Range("K71").Value = 60268.00199999999
$ws.Range("L71").Value = 118307.814
$ws.Range("M71").Value = -56524.00199999999
$ws.Range("N71").Value = -125795.814
$ws.Range("H74").Value = 12788.4
$ws.Range("J74").Value = 12788.4
$ws.Range("L74").Value = 12788.4
$ws.Range("N74").Value = -14536.4
$ws.Range("H77").Value = 12788.4
$ws.Range("J77").Value = 12788.4
$ws.Range("L77").Value = 38365.2
$ws.Range("N77").Value = -47101.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5714718.5
$ws.Range("I113").Value = 429.44446
$ws.Range("J113").Value = 20408604
$ws.Range("K113").Value = 1288.33338
$ws.Range("L113").Value = 61225812
$ws.Range("M113").Value = 881.66662
$ws.Range("N113").Value = -61230152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 9666.667
$ws.Range("I22").Value = 9666.667
$ws.Range("K22").Value = 9666.667
$ws.Range("M22").Value = -9137.667
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H82").Value = 47497.5
$ws.Range("I82").Value = 37000
$ws.Range("J82").Value = 57995
$ws.Range("K82").Value = 37000
$ws.Range("L82").Value = 57995
$ws.Range("M82").Value = -36617
$ws.Range("N82").Value = -58761
$ws.Range("H85").Value = 47497.5
$ws.Range("I85").Value = 37000
$ws.Range("J85").Value = 57995
$ws.Range("K85").Value = 37000
$ws.Range("L85").Value = 57995
$ws.Range("M85").Value = -35674
$ws.Range("N85").Value = -60647
$ws.Range("H113").Value = 9295.23
$ws.Range("I113").Value = 12447.333
$ws.Range("J113").Value = 2203
$ws.Range("K113").Value = 12447.333
$ws.Range("L113").Value = 2203
$ws.Range("M113").Value = -10277.333
$ws.Range("N113").Value = -6543
$ws.Range("H132").Value = 5892.6665
$ws.Range("I132").Value = 7853.4707
$ws.Range("J132").Value = 3328.5386
$ws.Range("K132").Value = 23560.4121
$ws.Range("L132").Value = 9985.6158
$ws.Range("M132").Value = -21030.4121
$ws.Range("N132").Value = -15045.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 590.44446
$ws.Range("J22").Value = 707.6923
$ws.Range("L22").Value = 707.6923
$ws.Range("N22").Value = -1297.6923
$ws.Range("H27").Value = 590.44446
$ws.Range("J27").Value = 707.6923
$ws.Range("L27").Value = 707.6923
$ws.Range("N27").Value = -921.6923
$ws.Range("H34").Value = 2000
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 27000
$ws.Range("J54").Value = 27000
$ws.Range("L54").Value = 27000
$ws.Range("N54").Value = -28040
